$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Units" (column B) for rows 30-37 from 0 to their delivered values
$ws.Range("B30").Value = 48
$ws.Range("B31").Value = 60
$ws.Range("B32").Value = 120
$ws.Range("B33").Value = 228
$ws.Range("B34").Value = 300
$ws.Range("B35").Value = 456
$ws.Range("B36").Value = 384
$ws.Range("B37").Value = 216

# Row 37 also had its Date corrected and its Status flipped from "Confirmed" to "Delivered"
$ws.Range("A37").Value = 44311
$ws.Range("D37").Value = "Delivered"
